# Update the EC (estado de cuenta) database: re-map the "Periodo Mora" values
# for the existing worker rows and refresh the "Salario Basico" amounts to the
# new value, as part 1 of the new estado de cuenta data load.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora column (E) re-ordering for rows 16-18:
#   E16: 2501 -> 2411
#   E17: 2412 -> 2412 (unchanged)
#   E18: 2411 -> 2501
$ws.Range("E16").Value = "2411"
$ws.Range("E17").Value = "2412"
$ws.Range("E18").Value = "2501"

# Salario Basico column (G) updated for rows 16-18: 1300000 -> 1423500
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
